$d = $word.ActiveDocument

# --- Change 1: Problem #1, part A. ---
# Remove the gramStart/gramEnd proofErr markers that bracket the lone "a"
# run and merge "a" + " man needs ... items?" back into a single run, while
# leaving the existing "A. " run untouched.
$p1 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "A. a man needs to get across a river on a boat. He has 3 things with him but the boat only fits himself and one other thing. What should he leave behind? What order should he take items?`r") {
        $p1 = $cand
    }
}
$rng1 = $p1.Range
# Exclude the trailing paragraph mark so only the runs (not the paragraph
# properties, e.g. the ListParagraph/numbering pPr) get replaced.
$sub1 = $d.Range($rng1.Start, $rng1.End - 1)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">A. </w:t></w:r><w:r><w:t>a man needs to get across a river on a boat. He has 3 things with him but the boat only fits himself and one other thing. What should he leave behind? What order should he take items?</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $sub1.InsertXML($xml1)

# --- Change 2: Problem #2, item 5 ---
# Locate the paragraph "5. " + bookmark + "A. The only guaranteed
# solution..." by its text (robust to any index drift) and split it into
# two paragraphs: "5. A. ..." and "    B. I actually ...", moving the
# _GoBack bookmark to trail the new "B." paragraph.
$p2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -eq "5. A. The only guaranteed solution is to pick out 20 socks to get the pairs you need to answer the questions.`r") {
        $p2 = $cand
    }
}
$rng2 = $p2.Range
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>5. A. The only guaranteed solution is to pick out 20 socks to get the pairs you need to answer the questions.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">    B. I actually  worked for a long time with my actual socks before the answer clicked that the only way to guarantee something decided at random is to pick  them all. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $rng2.InsertXML($xml2)
